$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.020.29'
$ws.Range('E2').Value = '  +3.07%  '

$ws.Range('D3').Value = '3.381.48'
$ws.Range('E3').Value = '  +1.73%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = "'594.11"
$ws.Range('E5').Value = '  +7.26%  '

$ws.Range('D6').Value = "'187.51"
$ws.Range('E6').Value = '  -0.20%  '

$ws.Range('D7').Value = "'0.602"
$ws.Range('E7').Value = '  +4.02%  '

$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').Value = "'0.185"
$ws.Range('E9').Value = '  +4.45%  '

$ws.Range('E10').Value = '  +2.04%  '

$ws.Range('D11').Value = "'47.68"
$ws.Range('E11').Value = '  +3.87%  '

$ws.Range('E12').Value = '  +5.78%  '

$ws.Range('D13').Value = '3.923.02'
$ws.Range('E13').Value = '  +1.91%  '

$ws.Range('D14').Value = "'639.92"
$ws.Range('E14').Value = '  +10.40%  '

$ws.Range('E15').Value = '  +2.15%  '

$ws.Range('D16').Value = '68.050.86'
$ws.Range('E16').Value = '  +3.20%  '

$ws.Range('D17').Value = '3.378.94'
$ws.Range('E17').Value = '  +1.78%  '

$ws.Range('E18').Value = '  +1.72%  '

$ws.Range('E19').Value = '  +1.89%  '

$ws.Range('D20').Value = "'11.12"
$ws.Range('E20').Value = '  +2.13%  '

$ws.Range('E21').Value = '  +2.29%  '

$ws.Range('E22').Value = '  -0.36%  '

$ws.Range('D23').Value = "'5.13"
$ws.Range('E23').Value = '  +2.54%  '

$ws.Range('D24').Value = "'99.77"
$ws.Range('E24').Value = '  +0.62%  '

$ws.Range('D25').Value = "'4.06"
$ws.Range('E25').Value = '  +2.93%  '

$ws.Range('E26').Value = '  +6.61%  '

$ws.Range('D27').Value = "'9.76"
$ws.Range('E27').Value = '  +4.66%  '

$ws.Range('D28').Value = "'32.90"
$ws.Range('E28').Value = '  +7.75%  '

$ws.Range('D29').Value = "'8.75"
$ws.Range('E29').Value = '  +4.34%  '

$ws.Range('E30').Value = '  +4.89%  '

$ws.Range('D31').Value = "'609.91"
$ws.Range('E31').Value = '  +6.16%  '

$ws.Range('D32').Value = "'3.83"
$ws.Range('E32').Value = '  +1.95%  '

$ws.Range('D33').Value = '4.037.15'
$ws.Range('E33').Value = '  +8.46%  '

$ws.Range('E34').Value = '  +2.61%  '

$ws.Range('E35').Value = '  +3.49%  '

$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = '  -0.05%  '

$ws.Range('D37').Value = "'56.35"
$ws.Range('E37').Value = '  +1.43%  '

$ws.Range('E38').Value = '  +6.49%  '

$ws.Range('E39').Value = '  +3.81%  '

$ws.Range('D40').Value = "'33.81"
$ws.Range('E40').Value = '  -0.60%  '

$ws.Range('E41').Value = '  +4.09%  '

$ws.Range('D42').Value = '0.0₃0706'
$ws.Range('E42').Value = '  +2.43%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = "'0.344"
$ws.Range('E43').Value = '  +2.93%  '

$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = "'3.41"
$ws.Range('E44').Value = '  +0.93%  '

$ws.Range('D45').Value = "'0.0424"
$ws.Range('E45').Value = '  +3.97%  '

$ws.Range('E46').Value = '  +1.71%  '

$ws.Range('D47').Value = "'2.61"
$ws.Range('E47').Value = '  +3.96%  '

$ws.Range('E48').Value = '  +12.76%  '

$ws.Range('E49').Value = '  +0.30%  '

$ws.Range('D50').Value = "'128.45"
$ws.Range('E50').Value = '  +1.40%  '

$ws.Range('D51').Value = "'7.78"
$ws.Range('E51').Value = '  +7.06%  '
